$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the inlineStr type in the target), otherwise Excel auto-converts
# them to numeric cells. We set NumberFormat to "@" (Text) before assignment,
# then ClearFormats() afterwards so no residual style index is left behind.
$numericLookingCells = @(
    "D5",
    "D6",
    "D9",
    "D13",
    "D20",
    "D21",
    "D24",
    "D25",
    "D27",
    "D28",
    "D31",
    "D33",
    "D36",
    "D37",
    "D40",
    "D41",
    "D43",
    "D46",
    "D48"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply all the new cell values
$ws.Range("D2").Value = "61.357.51"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "2.983.43"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "593.24"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "143.46"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "2.981.02"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "0.513"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  +5.01%  "
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").Value = "0.0000226"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "3.481.45"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "61.328.08"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "2.988.21"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "447.43"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").Value = "81.69"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "10.55"
$ws.Range("E25").Value = "  +6.31%  "
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").Value = "11.99"
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "7.17"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "27.15"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("D35").Value = "0.0₃0815"
$ws.Range("E35").Value = "  +3.78%  "
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "5.77"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").Value = "9.00"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "2.84"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").Value = "386.55"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "38.10"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.689.28"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").Value = "129.94"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E51").Value = "  -0.32%  "

# Remove the temporary text-format styling so these cells end up with no
# explicit style, consistent with the rest of the sheet
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).ClearFormats()
}
